$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 37 ("TCID 42" - LeaveBalance_48EmployeeCreation) is being split into four
# separate test-case rows (one per employee-creation batch). Insert three new
# rows right after the existing row so the block becomes rows 37-40, which
# pushes the old rows 38-44 down to 41-47.
$ws.Rows("38:40").Insert()

# The rows that used to follow (old TCID 43-49, now at rows 41-47) only need
# their TCID (column A) renumbered to keep the sequence contiguous; columns
# B:G already shifted down intact with the row insert above.
$ws.Range("A41").Value = "46"
$ws.Range("A42").Value = "47"
$ws.Range("A43").Value = "48"
$ws.Range("A44").Value = "49"
$ws.Range("A45").Value = "50"
$ws.Range("A46").Value = "51"
$ws.Range("A47").Value = "52"

# Row 37: first split (employees 1-25)
$ws.Range("A37").Value = "42"
$ws.Range("B37").Value = "LeaveBalance"
$ws.Range("C37").Value = "LeaveBalance"
$ws.Range("D37").Value = "com.darwinbox.leaves.Accural.Custom.LeaveBalance_48EmployeeCreation_1_25"
$ws.Range("E37").Value = "Accural//LeaveBalance.xlsx"
$ws.Range("F37").Value = "dummySheet"
$ws.Range("G37").Value = "All"

# Row 38: second split (employees 26-41)
$ws.Range("A38").Value = "43"
$ws.Range("B38").Value = "LeaveBalance"
$ws.Range("C38").Value = "LeaveBalance"
$ws.Range("D38").Value = "com.darwinbox.leaves.Accural.Custom.LeaveBalance_48EmployeeCreation_26_41"
$ws.Range("E38").Value = "Accural//LeaveBalance.xlsx"
$ws.Range("F38").Value = "dummySheet"
$ws.Range("G38").Value = "All"

# Row 39: third split (employees 83-100)
$ws.Range("A39").Value = "44"
$ws.Range("B39").Value = "LeaveBalance"
$ws.Range("C39").Value = "LeaveBalance"
$ws.Range("D39").Value = "com.darwinbox.leaves.Accural.Custom.LeaveBalance_48EmployeeCreation_83_100"
$ws.Range("E39").Value = "Accural//LeaveBalance.xlsx"
$ws.Range("F39").Value = "dummySheet"
$ws.Range("G39").Value = "All"

# Row 40: fourth split (employees 100-123)
$ws.Range("A40").Value = "45"
$ws.Range("B40").Value = "LeaveBalance"
$ws.Range("C40").Value = "LeaveBalance"
$ws.Range("D40").Value = "com.darwinbox.leaves.Accural.Custom.LeaveBalance_48EmployeeCreation_100_123"
$ws.Range("E40").Value = "Accural//LeaveBalance.xlsx"
$ws.Range("F40").Value = "dummySheet"
$ws.Range("G40").Value = "All"

# Update the active selection to match the saved view.
$ws.Range("A40").Select()
